$d = $word.ActiveDocument

# The bullet "A map of California cities with D. plexippus sightings..."
# currently renders the species name "D. plexippus" in plain (non-italic)
# text, unlike every other mention of "D. plexippus" in the document,
# which is italicized (scientific species names are conventionally set in
# italics). Bring this occurrence in line with the rest of the document by
# italicizing "D. " and "plexippus" (keeping the spell-check proofErr
# markers around "plexippus" untouched).

# Use an already-italicized "D. " / "plexippus" pair elsewhere in the same
# section of the document as the formatting source, so the new runs end up
# with both <w:i/> and <w:iCs/> (matching how the rest of the document
# marks up this species name), rather than just <w:i/>.
$srcRange = $d.Content
$srcRange.Find.Execute("D. plexippus might be") | Out-Null
$srcStart = $srcRange.Start
# $srcD is "D. " and $srcPlexippus is "plexippus"
$srcD = $d.Range($srcStart, $srcStart + 3)
$srcPlexippus = $d.Range($srcStart + 3, $srcStart + 12)

# Locate the target bullet point.
$bullet = $d.Content
$bullet.Find.Execute("A map of California cities with D. plexippus sightings") | Out-Null
$bulletStart = $bullet.Start
$bulletEnd = $bullet.End

# Within the bullet, find "D. " then "plexippus" right after it.
$afterPrefix = $d.Range($bulletStart, $bulletEnd)
$afterPrefix.Find.Execute("D. ") | Out-Null
$targetD = $d.Range($afterPrefix.Start, $afterPrefix.End)

$afterD = $d.Range($targetD.End, $bulletEnd)
$afterD.Find.Execute("plexippus") | Out-Null
$targetPlexippus = $d.Range($afterD.Start, $afterD.End)

# Transplant the italic (w:i + w:iCs) formatting onto the target text.
$targetD.FormattedText = $srcD.FormattedText
$targetPlexippus.FormattedText = $srcPlexippus.FormattedText
